$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 290, shifting existing rows 290-359 down to 291-360.
$ws.Rows.Item(290).EntireRow.Insert()

# Populate the newly inserted row 290 with the new weekly price record.
$ws.Cells.Item(290, 1).Value = 11
$ws.Cells.Item(290, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(290, 3).Value = "Bíobío"
$ws.Cells.Item(290, 4).Value = 45275
$ws.Cells.Item(290, 5).Value = 8
$ws.Cells.Item(290, 6).Value = 100112003
$ws.Cells.Item(290, 7).Value = "Ajo"
$ws.Cells.Item(290, 8).Value = "Chino"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 220
$ws.Cells.Item(290, 11).Value = 21000
$ws.Cells.Item(290, 12).Value = 22000
$ws.Cells.Item(290, 13).Value = 21545
$ws.Cells.Item(290, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(290, 15).Value = "China"
$ws.Cells.Item(290, 16).Value = 2154
$ws.Cells.Item(290, 17).Value = 10
$ws.Cells.Item(290, 18).Value = "Hortaliza"
